# Updated symbol list - apply new Price / Volume(1h) values
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Row 2
$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "328.16"
$ws.Range("E2").NumberFormat = "@"
$ws.Range("E2").Value = "-0.81%"
# Row 3
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "43.93"
$ws.Range("E3").NumberFormat = "@"
$ws.Range("E3").Value = "5.70%"
# Row 4
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = "5.475"
$ws.Range("E4").NumberFormat = "@"
$ws.Range("E4").Value = "-3.73%"
# Row 5
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "0.08086"
$ws.Range("E5").NumberFormat = "@"
$ws.Range("E5").Value = "-3.67%"
# Row 6
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "8.676"
$ws.Range("E6").NumberFormat = "@"
$ws.Range("E6").Value = "-1.56%"
# Row 7
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "4.291"
$ws.Range("E7").NumberFormat = "@"
$ws.Range("E7").Value = "-4.23%"
# Row 8
$ws.Range("E8").NumberFormat = "@"
$ws.Range("E8").Value = "-5.72%"
# Row 9
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "2.697"
$ws.Range("E9").NumberFormat = "@"
$ws.Range("E9").Value = "-8.20%"
# Row 10
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "0.9373"
$ws.Range("E10").NumberFormat = "@"
$ws.Range("E10").Value = "1.34%"
# Row 11
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.1191"
$ws.Range("E11").NumberFormat = "@"
$ws.Range("E11").Value = "-6.81%"
# Row 12
$ws.Range("E12").NumberFormat = "@"
$ws.Range("E12").Value = "-3.89%"
# Row 13
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "0.09661"
$ws.Range("E13").NumberFormat = "@"
$ws.Range("E13").Value = "2.28%"
# Row 14
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "0.04092"
$ws.Range("E14").NumberFormat = "@"
$ws.Range("E14").Value = "3.53%"
# Row 15
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "0.1068"
$ws.Range("E15").NumberFormat = "@"
$ws.Range("E15").Value = "0.53%"
# Row 16
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "0.001279"
$ws.Range("E16").NumberFormat = "@"
$ws.Range("E16").Value = "-1.57%"
# Row 17
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "0.005968"
$ws.Range("E17").NumberFormat = "@"
$ws.Range("E17").Value = "-2.32%"
# Row 18
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "3.568"
$ws.Range("E18").NumberFormat = "@"
$ws.Range("E18").Value = "4.23%"
# Row 20
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "8.619"
$ws.Range("E20").NumberFormat = "@"
$ws.Range("E20").Value = "-3.98%"
# Row 21
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "0.1348"
$ws.Range("E21").NumberFormat = "@"
$ws.Range("E21").Value = "-1.09%"
# Row 22
$ws.Range("E22").NumberFormat = "@"
$ws.Range("E22").Value = "-0.74%"
# Row 23
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "0.04362"
$ws.Range("E23").NumberFormat = "@"
$ws.Range("E23").Value = "-0.91%"
# Row 24
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "0.001235"
$ws.Range("E24").NumberFormat = "@"
$ws.Range("E24").Value = "-0.96%"
# Row 25
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "0.004305"
$ws.Range("E25").NumberFormat = "@"
$ws.Range("E25").Value = "-1.67%"
# Row 26
$ws.Range("E26").NumberFormat = "@"
$ws.Range("E26").Value = "3.44%"
# Row 27
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "0.0003999"
$ws.Range("E27").NumberFormat = "@"
$ws.Range("E27").Value = "0.15%"
# Row 39
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "0.02667"
$ws.Range("E39").NumberFormat = "@"
$ws.Range("E39").Value = "-5.83%"
# Row 40
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "0.05442"
$ws.Range("E40").NumberFormat = "@"
$ws.Range("E40").Value = "-1.54%"
# Row 41
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "0.007681"
$ws.Range("E41").NumberFormat = "@"
$ws.Range("E41").Value = "-3.41%"
# Row 42
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "0.009808"
$ws.Range("E42").NumberFormat = "@"
$ws.Range("E42").Value = "9.36%"
# Row 43
$ws.Range("E43").NumberFormat = "@"
$ws.Range("E43").Value = "-3.25%"
# Row 44
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "0.002091"
$ws.Range("E44").NumberFormat = "@"
$ws.Range("E44").Value = "-0.15%"
# Row 45
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "0.009903"
$ws.Range("E45").NumberFormat = "@"
$ws.Range("E45").Value = "-15.92%"
# Row 46
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "0.00007116"
$ws.Range("E46").NumberFormat = "@"
$ws.Range("E46").Value = "2.53%"
# Row 47
$ws.Range("E47").NumberFormat = "@"
$ws.Range("E47").Value = "0.09%"
# Row 48
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "0.003562"
$ws.Range("E48").NumberFormat = "@"
$ws.Range("E48").Value = "7.82%"
# Row 49
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "0.002276"
$ws.Range("E49").NumberFormat = "@"
$ws.Range("E49").Value = "-0.18%"
# Row 50
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "0.00002106"
$ws.Range("E50").NumberFormat = "@"
$ws.Range("E50").Value = "0.09%"
# Row 51
$ws.Range("E51").NumberFormat = "@"
$ws.Range("E51").Value = "0.09%"
